$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update J52's note to "N/A" (was "Seguire trabajando en fin de semana (rafael)")
$ws.Range("J52").Value = "N/A"

# Add new row 53 with the latest Roboflow annotation report entry (7/8/2025).
# Copy the previous data row's formatting down first so the new row matches
# the existing table styling (date format, borders, etc.).
$ws.Range("D52:J52").Copy($ws.Range("D53:J53"))
$ws.Range("D53").EntireRow.RowHeight = 15.6

$ws.Range("D53").Value = (Get-Date -Year 2025 -Month 8 -Day 7).Date
$ws.Range("E53").Value = 192
$ws.Range("F53").Value = 734
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 1012
$ws.Range("J53").Value = "N/A"

# Extend the Table1 list object to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("D4:J53"))

# Update view: scroll/selection moved
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 2
$ws.Range("E52").Select()
